# SwaadSutra_Daily_2026-01-19.xlsx update
# A new order (#19, Surekha Sonawane) came in. Insert it at the top of the
# "Daily Orders" log (row 2, pushing existing orders down by one row),
# add the matching line to "Items Breakdown" (Til Poli), and refresh the
# aggregate counters on "Summary".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Daily Orders: insert the new order as the first data row
# ---------------------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Daily Orders")
$wsOrders.Rows(2).Insert()

$wsOrders.Range("A2").Value = 19
$wsOrders.Range("B2").Value = "2026-01-19 09:18"
$wsOrders.Range("C2").Value = "Surekha Sonawane"
$wsOrders.Range("D2").Value = "A 808"

# Phone number / collection date are stored as text in this workbook, not
# as numbers / Excel dates - force text storage so they keep their literal
# digits (and the date string isn't reinterpreted as a date serial).
$wsOrders.Range("E2").NumberFormat = "@"
$wsOrders.Range("E2").Value = "935917349"

$wsOrders.Range("F2").Value = "Til Poli x4"
$wsOrders.Range("G2").Value = 120
$wsOrders.Range("H2").Value = "NEW"
$wsOrders.Range("I2").Value = "PENDING"

$wsOrders.Range("J2").NumberFormat = "@"
$wsOrders.Range("J2").Value = "2026-01-18"

$wsOrders.Range("K2").Value = "16:00"
$wsOrders.Range("L2").Value = ""
$wsOrders.Range("M2").Value = ""
$wsOrders.Range("N2").Value = ""

# ---------------------------------------------------------------
# Summary: bump Total Orders / New counts, add the new order's total
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A2").Value = 5
$wsSummary.Range("B2").Value = 5
$wsSummary.Range("G2").Value = 525

# ---------------------------------------------------------------
# Items Breakdown: add the "Til Poli" line (right after Wheat Chapati)
# ---------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("Items Breakdown")
$wsItems.Rows(3).Insert()
$wsItems.Range("A3").Value = "Til Poli"
$wsItems.Range("B3").Value = 4
$wsItems.Range("C3").Value = 120
